# Apply the "classFields" sheet reordering described in the commit diff.
# The underlying shared-string table was reordered (field name/type pairs for
# a few classes swapped places), which changes the text shown by several
# rows even though some of the row's other columns stay identical.
# Below we just set the resulting cell text directly; Excel will manage the
# shared-string table itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Row 2-4: SwaggerResourceConfig fields - "log" and "gatewayProperties" swap order
$ws.Range("B2").Value = "log"
$ws.Range("D2").Value = "org.slf4j.Logger"

$ws.Range("B4").Value = "gatewayProperties"
$ws.Range("D4").Value = "org.springframework.cloud.gateway.config.GatewayProperties"

# Row 9: AuthGlobalFilter.LOGGER field type
$ws.Range("D9").Value = "org.slf4j.Logger"

# Rows 11-14: ResourceServerConfig fields rotate order
$ws.Range("B11").Value = "restAuthenticationEntryPoint"
$ws.Range("D11").Value = "com.macro.mall.component.RestAuthenticationEntryPoint"

$ws.Range("B12").Value = "restfulAccessDeniedHandler"
$ws.Range("D12").Value = "com.macro.mall.component.RestfulAccessDeniedHandler"

$ws.Range("B13").Value = "ignoreUrlsConfig"
$ws.Range("D13").Value = "com.macro.mall.config.IgnoreUrlsConfig"

$ws.Range("B14").Value = "ignoreUrlsRemoveJwtFilter"
$ws.Range("D14").Value = "com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter"

# Rows 15-16: AuthorizationManager fields swap order
$ws.Range("B15").Value = "redisTemplate"
$ws.Range("D15").Value = "org.springframework.data.redis.core.RedisTemplate"

$ws.Range("B16").Value = "ignoreUrlsConfig"
$ws.Range("D16").Value = "com.macro.mall.config.IgnoreUrlsConfig"

# Row 17: IgnoreUrlsRemoveJwtFilter field name text
$ws.Range("B17").Value = "ignoreUrlsConfig"
